# "Base de datos Relacionada y diagrama terminado"
# Trabajamos Elmer Chandia / yobany Abrego
#
# The sheet had some leftover scratch/test values (hoy, Elmer, Yobany,
# Acetaminofen, Ibuprofeno, assorted numbers) in A25:C28 and C31:C34 that
# are cleaned up, and a missing field name (FKEstadoCivil) is filled in on
# the "EstadoCiviles" mini entity-relationship table at C13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the leftover scratch/testing data below the diagram tables.
$ws.Range("A25:C28").ClearContents()
$ws.Range("C31:C34").ClearContents()

# Add the missing foreign-key field name for the EstadoCiviles table.
$ws.Range("C13").Value = "FKEstadoCivil"

# Restore the view/selection to where the author left off working.
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 13
$ws.Range("O31").Select()
